$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.297.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.325.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.12"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.42%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.89%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.41%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.323.18"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.09%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.53"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.64%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.59"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.99%  "

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.741.13"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.15%  "

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.293.01"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000132"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.339.89"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.57"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.85"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.61%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.45%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.77"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.173"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.34%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.42%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.84%  "

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.46"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.48%  "

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "SuiNetwork"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.92%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.33%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.83%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.79%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +10.23%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.73%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.35%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.06%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "322.92"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +10.70%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.97"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.85"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.08%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0939"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.07"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.562"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.03%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.67%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0215"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +17.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.01"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.60%  "
